$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("580.25", "8.58", ...)
# must be pre-formatted as Text so Excel keeps them as the literal strings
# the source data uses (matches the original inlineStr cells), rather than
# silently converting them to numeric values.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D41", "D42", "D44", "D45", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.305.44"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.381.21"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "580.25"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "178.24"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Value = "0.197"
$ws.Range("E9").Value = "  +7.59%  "
$ws.Range("D10").Value = "0.586"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "48.33"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "0.0000283"
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").Value = "684.02"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "8.58"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "3.920.08"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "69.422.52"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "3.375.31"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "17.72"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "11.26"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "0.908"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "17.13"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "5.34"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "101.15"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "9.70"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "33.40"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "3.82"
$ws.Range("E31").Value = "  +16.31%  "
$ws.Range("D32").Value = "11.01"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "548.17"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "57.78"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "3.603.34"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +3.28%  "
$ws.Range("D39").Value = "35.24"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  +9.16%  "
$ws.Range("D41").Value = "3.32"
$ws.Range("E41").Value = "  +4.35%  "
$ws.Range("D42").Value = "2.70"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D44").Value = "0.0424"
$ws.Range("E44").Value = "  +3.05%  "
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "129.26"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "2.59"
$ws.Range("E51").Value = "  +0.63%  "
